$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 275
$wsExhibition.Range("F4").Value = 170
$wsExhibition.Range("F5").Value = 14

# Sheet "全部类型" (All Types) - fourth sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 275
$wsAll.Range("F5").Value = 14
